# This edit rearranges (permutes) the 30 data rows (rows 2-31) of the
# "Artfynd" sheet. The exact same set of records is kept - only the row
# each record lives in changes (a couple of rows keep their original
# position, e.g. row 30).
#
# Strategy:
#   1. Record, for every original data row, which columns have no cell at
#      all (as opposed to a cell that is simply blank) - this is needed
#      later so the rearranged rows end up with the same "shape" as their
#      original row.
#   2. Copy each source row (full A:AY range) into a temporary staging
#      area far below the used range. Range.Copy (rather than re-typing
#      values through .Value) preserves the original cell data types
#      exactly; this matters because several cells hold date-looking text
#      (e.g. "2023-06-21") that must remain plain text and not be
#      auto-converted into real dates.
#   3. Copy each staged row into its final destination row.
#   4. Re-apply the "no cell at all" information recorded in step 1, since
#      a rectangular Range.Copy always materializes every column in the
#      destination (even ones that had no cell in the source).
#   5. Clear the staging area so it doesn't remain in the saved file.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 31
$lastCol = "AY"
$lastColIdx = 51

# Maps destination row -> source row (both referring to the original,
# pre-edit positions of the data)
$mapping = @{
  2  = 9
  3  = 19
  4  = 21
  5  = 15
  6  = 4
  7  = 22
  8  = 7
  9  = 12
  10 = 13
  11 = 6
  12 = 17
  13 = 11
  14 = 3
  15 = 14
  16 = 24
  17 = 29
  18 = 16
  19 = 27
  20 = 10
  21 = 28
  22 = 23
  23 = 31
  24 = 20
  25 = 2
  26 = 25
  27 = 5
  28 = 8
  29 = 26
  30 = 30
  31 = 18
}

function Get-ColLetter([int]$colIndex) {
  return ($ws.Cells.Item(1, $colIndex).Address($false, $false, 1, $false) -replace '[0-9]+$', '')
}

# Step 1: figure out which columns have no cell at all in each original
# data row (a genuinely missing cell reads back as $null, whereas a cell
# that exists but is empty reads back as "").
$absentCols = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
  $cols = @()
  for ($c = 1; $c -le $lastColIdx; $c++) {
    if ($ws.Cells.Item($r, $c).Value() -eq $null) {
      $cols += (Get-ColLetter $c)
    }
  }
  $absentCols[$r] = $cols
}

$stagingOffset = 1000

# Step 2: copy every source row into a staging row far below the data.
# Clear the staging destination first so it starts out completely empty.
for ($r = $firstRow; $r -le $lastRow; $r++) {
  $src = $ws.Range("A" + $r + ":" + $lastCol + $r)
  $stageRow = $r + $stagingOffset
  $dst = $ws.Range("A" + $stageRow + ":" + $lastCol + $stageRow)
  $dst.Clear()
  $src.Copy($dst)
}

# Step 3: copy each staged row into its final destination row. Clear the
# destination first (Copy only overwrites cells that are populated in the
# source range; it does not blank out cells that are empty in the source
# but were non-empty at the destination before the copy).
foreach ($destRow in $mapping.Keys) {
  $sourceRow = $mapping[$destRow]
  $stageRow = $sourceRow + $stagingOffset
  $src = $ws.Range("A" + $stageRow + ":" + $lastCol + $stageRow)
  $dst = $ws.Range("A" + $destRow + ":" + $lastCol + $destRow)
  $dst.Clear()
  $src.Copy($dst)
}

# Step 4: a rectangular Range.Copy materializes a (blank) cell for every
# column of the destination, even the ones that had no cell at all in the
# source row. Remove those so each destination row again has exactly the
# same set of populated/blank-but-present/absent cells as its source row.
foreach ($destRow in $mapping.Keys) {
  $sourceRow = $mapping[$destRow]
  foreach ($col in $absentCols[$sourceRow]) {
    $ws.Range($col + $destRow).Clear()
  }
}

# Step 5: clear the staging area so it doesn't remain in the saved file.
$clearRange = $ws.Range("A" + ($stagingOffset + $firstRow) + ":" + $lastCol + ($stagingOffset + $lastRow))
$clearRange.Clear()
